# Insert a new "Prompts_EN" worksheet between "Data" and "Categorization".
# It joins the scenario (col I) and question (col J) text from the Data
# sheet into a single English prompt string, for rows 2-99.

$wb = $excel.ActiveWorkbook

$dataSheet = $wb.Worksheets.Item("Data")
$catSheet  = $wb.Worksheets.Item("Categorization")

# Adding "before" the Categorization sheet puts the new sheet right after
# Data, i.e. Data, Prompts_EN, Categorization.
$ws = $wb.Worksheets.Add($catSheet)
$ws.Name = "Prompts_EN"

$ws.Range("A1").Value = "Prompt"

$lastRow = $dataSheet.Cells.Item(1, 1).End(4).Row
if ($lastRow -lt 99) { $lastRow = 99 }

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Range("A$r").Formula = "=TEXTJOIN("" "", TRUE, Data!I${r}:J${r})"
}

# Make Prompts_EN the active sheet/tab with A2 selected, matching the
# author's saved view state.
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
